$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$t.Rows.Item($t.Rows.Count).Delete()
